# Added functionality to be able to check todays date against pay period date.
# Put today's (pay-period) date into E1 as a real date value/format, auto-size
# column E to fit it, and leave the selection on I11 (matching the recorded
# cursor position from the authoring session).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date value (serial 44711 == 2022-05-30) formatted as a short date
# (built-in numFmtId 14).
$ws.Range("E1").Value = 44711
$ws.Range("E1").NumberFormat = "mm-dd-yy"

# Auto-fit column E to the new content.
$ws.Columns("E:E").AutoFit()

# Restore the selection/active cell used when the sheet was saved.
$ws.Range("I11").Select() | Out-Null
